$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# Insert a new column before the "Settlement Remarks" column (AG) on the
# Settlement sheet, and give it a header.
$ws3.Columns("AG:AG").Insert()
$ws3.Range("AG1").Value = "Delete Settlement textbox"

# Duplicate row 3 into a brand-new row 4 (QA_TestCase_Auto_Optimus_3_1_3),
# copying formats first and then values so that blank cells keep their style.
$src = $ws3.Range("A3:AK3")
$dst = $ws3.Range("A4:AK4")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$src.Copy()
$dst.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$ws3.Range("A4").Value = "QA_TestCase_Auto_Optimus_3_1_3"
$ws3.Range("AG4").Value = "Delete Settlement"

# Update the active selections / active sheet to match the final state.
[void]$ws2.Range("B7").Select()
[void]$ws3.Activate()
[void]$ws3.Range("B6").Select()
